$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange

$full = $tr.Text
$anchor = $full.IndexOf("мероприиятиях")

# Original run boundaries (0-based offsets into $full), left to right:
#   run A: "мероприиятиях"   length 13
#   run B: " и "             length 3
#   run C: "выводид"         length 7
#   run D: " их  на экран."  length 14
# New contents for the same four runs (same run count/order, text re-wrapped):
#   run A: "мероприятиях "
#   run B: "и "
#   run C: "выводит "
#   run D: "их  на экран."

$offA = $anchor
$offB = $anchor + 13
$offC = $anchor + 13 + 3
$offD = $anchor + 13 + 3 + 7

# Apply right-to-left so earlier offsets stay valid while lengths change.
$tr.Characters($offD + 1, 14).Text = "их  на экран."
$tr.Characters($offC + 1, 7).Text  = "выводит "
$tr.Characters($offB + 1, 3).Text  = "и "
$tr.Characters($offA + 1, 13).Text = "мероприятиях "
